$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.688.10"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "'1.889.34"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -1.29%  "
$ws.Range("D5").Value = "'313.11"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("D7").Value = "'0.4845"
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("D8").Value = "'0.3789"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("D9").Value = "'0.07328"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").Value = "'0.9172"
$ws.Range("E10").Value = "  -1.98%  "
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("D12").Value = "'0.07690"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("D13").Value = "'1.906.70"
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").Value = "'5.458"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "'6.594"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").Value = "'90.88"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("D18").Value = "'0.000008786"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").Value = "'27.719.47"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'14.50"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("D22").Value = "'5.114"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'2.133.11"
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").Value = "'1.913"
$ws.Range("E25").Value = "  -1.68%  "
$ws.Range("D26").Value = "'153.24"
$ws.Range("E26").Value = "  -2.04%  "
$ws.Range("D27").Value = "'18.36"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "'2.116"
$ws.Range("E28").Value = "  +2.83%  "
$ws.Range("D29").Value = "'115.75"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "'4.891"
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("D31").Value = "'0.08927"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").Value = "'3.148"
$ws.Range("E32").Value = "  -5.38%  "
$ws.Range("D33").Value = "'1.219"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").Value = "'0.7608"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("D35").Value = "'4.627"
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("D36").Value = "'0.02033"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "'2.537"
$ws.Range("E37").Value = "  -6.27%  "
$ws.Range("D38").Value = "'1.091"
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("D39").Value = "'0.05253"
$ws.Range("E39").Value = "  -2.45%  "
$ws.Range("E40").Value = "  -3.30%  "
$ws.Range("D41").Value = "'2.976"
$ws.Range("D42").Value = "'6.921"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("D43").Value = "'0.1519"
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("D44").Value = "'8.312"
$ws.Range("E44").Value = "  -2.75%  "
$ws.Range("D45").Value = "'109.62"
$ws.Range("E45").Value = "  +4.52%  "
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("D47").Value = "'0.4775"
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("D49").Value = "'1.632"
$ws.Range("E49").Value = "  -2.20%  "
$ws.Range("D50").Value = "'67.28"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").Value = "'0.06055"
$ws.Range("E51").Value = "  -0.81%  "
